$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column N to the new column O for the relevant rows
$ws.Range("N3:N14").Copy()
$ws.Range("O3:O14").PasteSpecial(-4122)

# Set the new values in column O
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 70.636215334420882
$ws.Range("O6").Value = 107.1
$ws.Range("O7").Value = 55.452054794520542
$ws.Range("O8").Value = 84.375
$ws.Range("O9").Value = 120.48192771084337
$ws.Range("O10").Value = 109.53346855983774
$ws.Range("O11").Value = 147.7690288713911
$ws.Range("O12").Value = 25.545675020210183
$ws.Range("O13").Value = 82.457854874175425
$ws.Range("O14").Value = 15.384615384615385

$excel.CutCopyMode = 0
